$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Price/Volume columns to text format so numeric-looking strings
# (e.g. "5.23", "0.0625") are not auto-converted to numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "36.992.12"
$ws.Range("E2").Value = "  +1.20%  "
$ws.Range("D3").Value = "2.049.63"
$ws.Range("E3").Value = "  -2.42%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "248.87"
$ws.Range("E5").Value = "  -1.51%  "
$ws.Range("D6").Value = "0.668"
$ws.Range("E6").Value = "  +1.67%  "
$ws.Range("D8").Value = "55.34"
$ws.Range("E8").Value = "  +16.83%  "
$ws.Range("D9").Value = "60.63"
$ws.Range("E9").Value = "  +0.71%  "
$ws.Range("D10").Value = "0.382"
$ws.Range("E10").Value = "  +1.26%  "
$ws.Range("E11").Value = "  +4.94%  "
$ws.Range("E12").Value = "  +6.04%  "
$ws.Range("D13").Value = "14.95"
$ws.Range("E13").Value = "  +2.01%  "
$ws.Range("D14").Value = "2.352.46"
$ws.Range("E14").Value = "  -2.17%  "
$ws.Range("E15").Value = "  -1.47%  "
$ws.Range("D16").Value = "5.23"
$ws.Range("E16").Value = "  +2.91%  "
$ws.Range("D17").Value = "2.050.21"
$ws.Range("E17").Value = "  -2.30%  "
$ws.Range("D18").Value = "36.941.59"
$ws.Range("E18").Value = "  +1.27%  "
$ws.Range("D19").Value = "0.0₃0920"
$ws.Range("E19").Value = "  +11.10%  "
$ws.Range("D20").Value = "73.00"
$ws.Range("E20").Value = "  +0.37%  "
$ws.Range("D21").Value = "14.21"
$ws.Range("E21").Value = "  +7.95%  "
$ws.Range("D22").Value = "5.33"
$ws.Range("E22").Value = "  +3.34%  "
$ws.Range("D23").Value = "236.02"
$ws.Range("E23").Value = "  -1.58%  "
$ws.Range("E24").Value = "  -0.04%  "
$ws.Range("D25").Value = "2.41"
$ws.Range("E25").Value = "  -2.10%  "
$ws.Range("D26").Value = "170.14"
$ws.Range("E26").Value = "  -0.37%  "
$ws.Range("D27").Value = "8.95"
$ws.Range("E27").Value = "  -1.93%  "
$ws.Range("D28").Value = "20.01"
$ws.Range("E28").Value = "  -6.98%  "
$ws.Range("D29").Value = "1.97"
$ws.Range("E29").Value = "  -0.80%  "
$ws.Range("D30").Value = "0.125"
$ws.Range("E30").Value = "  +1.25%  "
$ws.Range("D31").Value = "4.59"
$ws.Range("E31").Value = "  +3.10%  "
$ws.Range("D32").Value = "0.0625"
$ws.Range("E32").Value = "  +1.72%  "
$ws.Range("D33").Value = "1.05"
$ws.Range("E33").Value = "  +5.47%  "
$ws.Range("D34").Value = "4.36"
$ws.Range("E34").Value = "  +6.92%  "
$ws.Range("E35").Value = "  -0.14%  "
$ws.Range("D36").Value = "0.0872"
$ws.Range("E36").Value = "  -5.45%  "
$ws.Range("E37").Value = "  -6.18%  "
$ws.Range("E38").Value = "  -4.83%  "
$ws.Range("D39").Value = "1.35"
$ws.Range("E39").Value = "  +0.53%  "
$ws.Range("E40").Value = "  +22.22%  "
$ws.Range("D41").Value = "17.83"
$ws.Range("E41").Value = "  +12.06%  "
$ws.Range("D42").Value = "0.0224"
$ws.Range("E42").Value = "  +0.65%  "
$ws.Range("E43").Value = "  -2.51%  "
$ws.Range("D44").Value = "96.41"
$ws.Range("E44").Value = "  -1.32%  "
$ws.Range("D46").Value = "4.13"
$ws.Range("E46").Value = "  +46.61%  "
$ws.Range("D47").Value = "13.52"
$ws.Range("E47").Value = "  -53.12%  "
$ws.Range("D48").Value = "2.38"
$ws.Range("E48").Value = "  +6.28%  "
$ws.Range("D49").Value = "1.292.93"
$ws.Range("E49").Value = "  -2.53%  "
$ws.Range("D50").Value = "2.91"
$ws.Range("E50").Value = "  +2.05%  "
$ws.Range("E51").Value = "  +6.81%  "

# Restore the original (default) cell style now that the text values are set.
$ws.Range("D2:E51").Style = "Normal"
